$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.329.36'
$ws.Range("E2").Value = '  +1.09%  '
$ws.Range("D3").Value = '1.920.87'
$ws.Range("E3").Value = '  +0.59%  '
$ws.Range("D4").Value = '''1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = '''0.8096'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.49%  '
$ws.Range("D6").Value = '''244.49'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.15%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").Value = '''0.3243'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.83%  '
$ws.Range("D9").Value = '''27.19'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.08%  '
$ws.Range("D10").Value = '''0.07105'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.21%  '
$ws.Range("D11").Value = '''0.7845'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.94%  '
$ws.Range("D12").Value = '''0.08099'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.75%  '
$ws.Range("D13").Value = '1.923.16'
$ws.Range("E13").Value = '  +0.84%  '
$ws.Range("D14").Value = '''5.416'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.18%  '
$ws.Range("D15").Value = '''94.94'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.32%  '
$ws.Range("D16").Value = '30.316.59'
$ws.Range("E16").Value = '  +1.11%  '
$ws.Range("D17").Value = '''14.33'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.05%  '
$ws.Range("D18").Value = '''6.029'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.84%  '
$ws.Range("D19").Value = '''249.59'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.78%  '
$ws.Range("D20").Value = '''0.000007828'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.28%  '
$ws.Range("D21").Value = '2.178.16'
$ws.Range("E21").Value = '  +0.99%  '
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("E23").Value = '  +0.22%  '
$ws.Range("D24").Value = '''7.898'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +15.30%  '
$ws.Range("D25").Value = '''0.1630'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +18.24%  '
$ws.Range("D26").Value = '''9.520'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.83%  '
$ws.Range("E27").Value = '  -0.53%  '
$ws.Range("D28").Value = '''19.11'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.33%  '
$ws.Range("D29").Value = '''2.128'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.85%  '
$ws.Range("D30").Value = '''1.374'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.27%  '
$ws.Range("D31").Value = '''1.540'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.43%  '
$ws.Range("D32").Value = '''4.366'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.41%  '
$ws.Range("D33").Value = '''0.05648'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.45%  '
$ws.Range("D34").Value = '''4.135'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.43%  '
$ws.Range("D35").Value = '''1.306'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.23%  '
$ws.Range("D36").Value = '''0.7423'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.48%  '
$ws.Range("D37").Value = '''1.000'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.03%  '
$ws.Range("D38").Value = '''2.720'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.46%  '
$ws.Range("D39").Value = '''0.01951'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.10%  '
$ws.Range("D40").Value = '''2.819'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.04%  '
$ws.Range("D41").Value = '''0.4481'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.55%  '
$ws.Range("D42").Value = '''73.83'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.88%  '
$ws.Range("D43").Value = '''5.985'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.18%  '
$ws.Range("D44").Value = '''0.8546'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.78%  '
$ws.Range("D45").Value = '''1.934'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.48%  '
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").Value = '''1.002'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.11%  '
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '1.039.05'
$ws.Range("E47").Value = '  +6.01%  '
$ws.Range("D48").Value = '''103.16'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.83%  '
$ws.Range("D49").Value = '''9.974'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.55%  '
$ws.Range("D50").Value = '''7.648'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.76%  '
$ws.Range("D51").Value = '2.067.90'
$ws.Range("E51").Value = '  +0.53%  '
